# Disponibilidad.xlsx automatic update — 02-05-2021 05-15-33
# Mirrors: update the timestamp on the last existing 14-row block
# (rows 352:365) and append a brand-new 14-row block (rows 366:379)
# that repeats the same Nombre/URL/Disponibilidad pattern with a
# newer timestamp, wiring up the matching hyperlinks in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Existing rows 352:365 — only the Fecha (column D) timestamp
#    changes, in place, for every row of the last block.
# ---------------------------------------------------------------
$ws.Range("D352:D365").Value() = 44232.19803153935

# ---------------------------------------------------------------
# 2) Append rows 366:379 — a new copy of the same 14-row pattern.
# ---------------------------------------------------------------
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
# hyperlink target address (no trailing #/ fragment — that part is
# carried as the hyperlink's SubAddress/location instead)
$linkAddresses = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$startRow = 366
for ($i = 0; $i -lt 14; $i++) {
    $row = $startRow + $i

    $ws.Range("A$row").Value() = $names[$i]
    $ws.Range("B$row").Value() = $urls[$i]
    $ws.Range("C$row").Value() = "Disponible"
    $ws.Range("D$row").Value() = 44232.21908939274
    $ws.Range("D$row").NumberFormat() = "YYYY-MM-DD HH:MM:SS"
}

# ---------------------------------------------------------------
# 3) Wire up the column-B hyperlinks for the new rows. Row 374 is
#    the "MapStore" row whose display text carries a trailing
#    "#/" fragment — that fragment is the hyperlink SubAddress.
# ---------------------------------------------------------------
for ($i = 0; $i -lt 14; $i++) {
    $row = $startRow + $i
    if ($row -eq 374) {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $linkAddresses[$i], "/")
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $linkAddresses[$i])
    }
    # Hyperlinks.Add re-styles the cell with a fresh ad-hoc style;
    # put it back on the shared "Hyperlink" cell style used by
    # every other link cell in the column.
    $ws.Range("B$row").Style() = "Hyperlink"
}
